# Adding the change-password option to the admin and unite the select
# role of the user in the creation process.
#
# Users sheet columns: A name | B mador | C id | D password | E type |
#                       F is_active | G is_admin | H last_login | I profile_image

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (שחר יהודה יאני) ---------------------------------------------
$ws.Cells.Item(2, 5).Value = 0                        # E2 type
$ws.Cells.Item(2, 8).Value = "23:24 28/12/2025"       # H2 last_login

# --- Row 3 (רוני) -------------------------------------------------------
$ws.Cells.Item(3, 2).Value = "M"                      # B3 mador: M1 -> M

# is_active / is_admin are now real booleans instead of 0/1 numbers
$ws.Cells.Item(3, 6).Value = $false                   # F3 is_active
$ws.Cells.Item(3, 7).Value = $false                   # G3 is_admin

# last_login becomes plain text instead of a formatted date serial, so
# drop the old number-format style before writing the new value
$ws.Cells.Item(3, 8).ClearFormats()
$ws.Cells.Item(3, 8).Value = "23:45 28/12/2025"       # H3 last_login

# profile_image column removed for this user
$ws.Cells.Item(3, 9).ClearContents()                  # I3 profile_image

# --- Row 4 (יובל) --------------------------------------------------------
$ws.Cells.Item(4, 5).Value = 0                        # E4 type

$ws.Cells.Item(4, 6).Value = $false                   # F4 is_active
$ws.Cells.Item(4, 7).Value = $true                    # G4 is_admin

$ws.Cells.Item(4, 8).ClearFormats()
$ws.Cells.Item(4, 8).Value = "23:45 28/12/2025"       # H4 last_login

$ws.Cells.Item(4, 9).ClearContents()                  # I4 profile_image

# --- Row 5 (איתי) - new user row added -----------------------------------
$ws.Cells.Item(5, 1).Value = "איתי"                   # A5 name
$ws.Cells.Item(5, 2).Value = "M1"                     # B5 mador
$ws.Cells.Item(5, 3).Value = 4                        # C5 id

# password must stay textual ("123"), not be auto-coerced to a number
$ws.Cells.Item(5, 4).Value = "'123"                   # D5 password
$ws.Cells.Item(5, 4).ClearFormats()

$ws.Cells.Item(5, 5).Value = 2                        # E5 type
$ws.Cells.Item(5, 6).Value = $false                   # F5 is_active
$ws.Cells.Item(5, 7).Value = $false                   # G5 is_admin
$ws.Cells.Item(5, 8).Value = "23:46 28/12/2025"       # H5 last_login
# I5 profile_image intentionally left empty
